# Edits to micromorphological analysis and graphics to reflect splitting of
# forest and river habitats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a few data-entry typos -------------------------------------------------

# Row 70 (Londono_and_Quintero_214): Photographer_leaves "EKm" -> "EKM"
$ws.Cells.Item(70, 11).Value = "EKM"

# Row 86 (Pohl_and_Clark_13930): Floral_bracts_imaged "M" -> "No"
$ws.Cells.Item(86, 10).Value = "No"

# Row 145 (Clark_and_Oliveira_913): Floral_bracts_imaged "x" -> "No"
$ws.Cells.Item(145, 10).Value = "No"

# --- Insert a new specimen row --------------------------------------------------
# New record "Nelson_6026" belongs (alphabetically, within the G_weberbaueri
# group) right before the existing row 182 ("Seibert_2068"), so insert a row
# there and shift the remaining rows down.

$ws.Rows.Item(182).Insert()

$ws.Cells.Item(182, 1).Value = "Nelson_6026"
$ws.Cells.Item(182, 2).Value = "G_weberbaueri"
$ws.Cells.Item(182, 3).Value = "x"
$ws.Cells.Item(182, 4).Value = "Brazil"
$ws.Cells.Item(182, 5).Value = "Eastern_South_America"
$ws.Cells.Item(182, 6).Value = "Leaning_climbing"
$ws.Cells.Item(182, 7).Value = "Forest"
$ws.Cells.Item(182, 8).Value = "US"
$ws.Cells.Item(182, 9).Value = "Yes"
$ws.Cells.Item(182, 10).Value = "Yes"
$ws.Cells.Item(182, 11).Value = "EKM"
$ws.Cells.Item(182, 12).Value = "EKM"
$ws.Cells.Item(182, 13).Value = "Pt"
$ws.Cells.Item(182, 14).Value = "Pt"
$ws.Cells.Item(182, 15).Value = "EKM"
$ws.Cells.Item(182, 16).Value = "EKM"

# --- Update view state (frozen pane / active selection) -------------------------
$ws.Application.ActiveWindow.ScrollRow = 123
$ws.Range("I134").Select()
